# Add a new worksheet "NEGOCIAÇÕES Descrição" after the last existing sheet
# (Plan1, Plan2, Plan3) and populate it with a small header/value pair.
# This reproduces the commit that exposes a SheetCell naming bug with
# unicode characters in the sheet name.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "NEGOCIAÇÕES Descrição"

$ws.Range("A1").Value = "Negociações"
$ws.Range("B1").Value = 10

# Match the selection left behind on the new sheet in the target workbook.
$ws.Range("B2").Select()

# Keep the original first sheet as the active/selected tab, as in the
# target workbook (tabSelected stays on "Plan1").
$wb.Worksheets.Item("Plan1").Activate()
